$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1: drop the ToDo5 row (row 6) so the sheet ends at row 5.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 2. Sheet2: drop rows 4-7 (ToDo5..ToDo8) so the sheet ends at row 3.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows("4:7").Delete()

# ---------------------------------------------------------------------------
# 3. Add a brand-new Sheet3 right after Sheet2 with the TodoBackend links.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

$ws3.Columns.Item(1).ColumnWidth = 13.83
$ws3.Columns.Item(2).ColumnWidth = 64.16666666666667

$ws3.Range("A1").Value = "Name"
$ws3.Range("B1").Value = "baseURL"

$ws3.Range("A2").Value = "Completed"
$ws3.Range("A3").Value = "Active"
$ws3.Range("A4").Value = "All"

$ws3.Range("B2").Value = "https://www.todobackend.com/client/index.html?https://mysterious-thicket-31854.herokuapp.com/#/completed"
$ws3.Range("B3").Value = "https://www.todobackend.com/client/index.html?https://mysterious-thicket-31854.herokuapp.com/#/active"
$ws3.Range("B4").Value = "https://www.todobackend.com/client/index.html?https://mysterious-thicket-31854.herokuapp.com/#/"

# Thin box border around the whole A1:B4 table.
$ws3.Range("A1:B4").Borders.LineStyle = 1

# Live hyperlinks for the three URL cells (base target + in-page fragment).
$baseUrl = "https://www.todobackend.com/client/index.html?https://mysterious-thicket-31854.herokuapp.com/"
$ws3.Hyperlinks.Add($ws3.Range("B2"), $baseUrl, "/completed")
$ws3.Hyperlinks.Add($ws3.Range("B3"), $baseUrl, "/active")
$ws3.Hyperlinks.Add($ws3.Range("B4"), $baseUrl, "/")

$ws3.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. View/selection bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------
$ws2.Range("D10").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
